$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '25.843.27'
$ws.Range('E2').Value = '  -1.25%  '

# Row 3
$ws.Range('D3').Value = '1.636.88'
$ws.Range('E3').Value = '  -1.38%  '

# Row 4
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.42'
$ws.Range('E5').Value = '  -1.37%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5021'
$ws.Range('E6').Value = '  -2.77%  '

# Row 7
$ws.Range('E7').Value = '  -0.21%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2570'
$ws.Range('E8').Value = '  -0.82%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06418'

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.61'
$ws.Range('E10').Value = '  -1.69%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07687'
$ws.Range('E11').Value = '  -1.35%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.245'
$ws.Range('E12').Value = '  -1.23%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.635.09'
$ws.Range('E13').Value = '  -1.47%  '

# Row 14
$ws.Range('D14').Value = '1.860.84'
$ws.Range('E14').Value = '  -1.43%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5462'
$ws.Range('E15').Value = '  -1.73%  '

# Row 16
$ws.Range('E16').Value = '  -1.86%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '63.54'
$ws.Range('E17').Value = '  -1.14%  '

# Row 18
$ws.Range('D18').Value = '25.861.15'
$ws.Range('E18').Value = '  -1.32%  '

# Row 19
$ws.Range('E19').Value = '  -0.26%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '203.33'
$ws.Range('E20').Value = '  -4.15%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.303'
$ws.Range('E21').Value = '  -2.63%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.955'
$ws.Range('E22').Value = '  -0.91%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.989'
$ws.Range('E23').Value = '  +0.37%  '

# Row 24
$ws.Range('E24').Value = '  -0.20%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.936'
$ws.Range('E25').Value = '  +10.18%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '141.42'
$ws.Range('E26').Value = '  -2.04%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1146'
$ws.Range('E27').Value = '  -1.59%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.67'
$ws.Range('E28').Value = '  -1.06%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.709'
$ws.Range('E29').Value = '  -3.88%  '

# Row 30
$ws.Range('E30').Value = '  -1.39%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.04962'
$ws.Range('E31').Value = '  -5.97%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.270'
$ws.Range('E32').Value = '  -2.87%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.189'
$ws.Range('E33').Value = '  -0.99%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.530'
$ws.Range('E34').Value = '  -2.64%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.352'
$ws.Range('E35').Value = '  -0.78%  '

# Row 36
$ws.Range('D36').Value = '1.177.17'
$ws.Range('E36').Value = '  +0.61%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.8937'
$ws.Range('E37').Value = '  -3.71%  '

# Row 38
$ws.Range('E38').Value = '  -5.10%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5566'
$ws.Range('E39').Value = '  -1.78%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01558'
$ws.Range('E40').Value = '  -2.23%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.557'
$ws.Range('E41').Value = '  -0.22%  '

# Row 42
$ws.Range('E42').Value = '  -0.26%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.638'
$ws.Range('E43').Value = '  -1.01%  '

# Row 44
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '99.55'
$ws.Range('E44').Value = '  -0.97%  '

# Row 45
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8043'
$ws.Range('E45').Value = '  -4.72%  '

# Row 46
$ws.Range('D46').Value = '1.773.31'
$ws.Range('E46').Value = '  -1.34%  '

# Row 47
$ws.Range('E47').Value = '  -2.37%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4510'
$ws.Range('E48').Value = '  -0.55%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.005'
$ws.Range('E49').Value = '  +0.00%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '54.84'
$ws.Range('E50').Value = '  -1.85%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05038'
$ws.Range('E51').Value = '  -0.33%  '

